# Master_Bug_List.xlsx edit:
# docs(report): correct BUG_PROD_001 from invalid XSS to UI/UX missing keyword
# issue and update severity metrics
#
# BUG_PROD_001 (row 5) was originally reported as an XSS vulnerability on the
# product Search bar. QA re-verified the bug: it is not an XSS hole (the HTML
# is not executed) but rather a UI/UX defect where the "SEARCHED PRODUCTS"
# results header never echoes back the keyword the user typed. This script
# rewrites the row's content to describe the corrected bug and lowers the
# severity/priority metrics accordingly, then tidies up the now-irrelevant
# conditional-formatting rule and view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 / BUG_PROD_001: rewrite the bug report body -------------------
# A5 (Bug ID) and I5 (Status) are unchanged ("BUG_PROD_001" / "Open").

$ws.Range("B5").Value = "TC_PROD_001"

$ws.Range("C5").Value = "[Cập nhật] Lỗi UX: Không hiển thị lại từ khóa đã tìm kiếm trên màn hình kết quả"

$ws.Range("D5").Value = "Major"
$ws.Range("E5").Value = "Medium"

$ws.Range("F5").Value = "1. Vào trang Products`n2. Nhập từ khóa bất kỳ vào ô Search`n3. Nhấn Tìm kiếm"

$ws.Range("G5").Value = "Tiêu đề chỉ hiển thị chung chung là ""SEARCHED PRODUCTS"", không chứa từ khóa người dùng đã nhập."

$ws.Range("H5").Value = "Phải hiển thị rõ từ khóa để người dùng xác nhận, VD: ""SEARCHED PRODUCTS FOR 'Blue Top'""."

# The corrected description wraps onto more lines, so the row grows taller.
$ws.Rows.Item(5).RowHeight = 61.5

# --- Conditional formatting clean-up --------------------------------------
# Drop the now-redundant single-cell "Major" highlight rule that only ever
# applied to D2; the blanket D2:D11 rule set (Major/Minor/"cri") already
# covers every severity cell, including D2, so this extra rule is removed.
$fcs = $ws.Range("D2").FormatConditions
for ($i = $fcs.Count; $i -ge 1; $i--) {
    $rule = $fcs.Item($i)
    if ($rule.AppliesTo.Address() -eq "`$D`$2") {
        $rule.Delete()
    }
}

# --- View state -------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 2
$win.TopLeftCell = $ws.Range("C2")
$ws.Range("I5").Select()

Write-Host "BUG_PROD_001 rewritten as UX keyword-echo defect; severity/priority downgraded."
